$wb = $excel.ActiveWorkbook

# Rename the first sheet from "Sheet1" to "Sort by Tier"
$sheet1 = $wb.Worksheets.Item(1)
$sheet1.Name = "Sort by Tier"

# Make "Sort by Tier" the active/tab-selected sheet with B29 selected.
# (This also causes "Sort by ID" to lose its previous tabSelected flag,
# while it keeps its own D9 selection.)
$sheet1.Activate()
$sheet1.Range("B29").Select()
